$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $row = 44

    # A44 - Trade #
    $ws.Cells.Item($row, 1).Value = 43

    # B44 - Date (force text so "2026-02-17" isn't auto-converted to a date serial)
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).ClearFormats()

    # C44 - Time
    $ws.Cells.Item($row, 3).Value = "08:33:38"

    # D44 - Strategy
    $ws.Cells.Item($row, 4).Value = "MarketMaking"

    # E44 - Side
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # F44 - Entry Price
    $ws.Cells.Item($row, 6).Value = 0.46

    # G44 - Exit Price (blank / open trade, stored as empty text)
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).ClearFormats()

    # H44 - Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # I44 - P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # J44 - P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # K44 - Capital After
    $ws.Cells.Item($row, 11).Value = 99.77598934440597

    # L44 - Entry Slippage (bps)
    $ws.Cells.Item($row, 12).Value = 0

    # M44 - Exit Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0

    # N44 - Confidence
    $ws.Cells.Item($row, 14).Value = 0.6

    # O44 - Entry Reason
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # P44 - Exit Reason (blank / open trade, stored as empty text)
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).ClearFormats()

    # Q44 - Duration (min)
    $ws.Cells.Item($row, 17).Value = 0
}
